$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.234.59'
$ws.Range("E2").Value = '  -4.84%  '
$ws.Range("D3").Value = '3.316.16'
$ws.Range("E3").Value = '  -5.17%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.01'
$ws.Range("E5").Value = '  -3.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.20'
$ws.Range("E6").Value = '  -5.27%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.312.92'
$ws.Range("E8").Value = '  -5.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.476'
$ws.Range("E9").Value = '  -2.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.20'
$ws.Range("E10").Value = '  -5.40%  '
$ws.Range("E11").Value = '  -5.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.374'
$ws.Range("E12").Value = '  -4.23%  '
$ws.Range("D13").Value = '3.880.33'
$ws.Range("E13").Value = '  -5.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.119'
$ws.Range("E14").Value = '  -1.29%  '
$ws.Range("D15").Value = '3.318.11'
$ws.Range("E15").Value = '  -5.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000169'
$ws.Range("E16").Value = '  -6.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.82'
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("D18").Value = '61.325.98'
$ws.Range("E18").Value = '  -4.70%  '
$ws.Range("E19").Value = '  -9.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.58'
$ws.Range("E20").Value = '  -3.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.17'
$ws.Range("E21").Value = '  -2.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '355.65'
$ws.Range("E22").Value = '  -8.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.554'
$ws.Range("E23").Value = '  -4.79%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = '3.446.28'
$ws.Range("E25").Value = '  -5.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.22'
$ws.Range("E26").Value = '  -5.48%  '
$ws.Range("E27").Value = '  -7.37%  '
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.20'
$ws.Range("E29").Value = '  -2.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.47'
$ws.Range("E30").Value = '  -1.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  -3.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.11'
$ws.Range("E32").Value = '  -6.91%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.148'
$ws.Range("E34").Value = '  -5.09%  '
$ws.Range("D35").Value = '3.341.06'
$ws.Range("E35").Value = '  -5.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.47'
$ws.Range("E36").Value = '  +2.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.35'
$ws.Range("E37").Value = '  -4.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.73'
$ws.Range("E38").Value = '  -2.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '162.50'
$ws.Range("E39").Value = '  -1.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.50'
$ws.Range("E40").Value = '  -3.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0756'
$ws.Range("E41").Value = '  -4.02%  '
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.84'
$ws.Range("E43").Value = '  -2.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.748'
$ws.Range("E44").Value = '  -7.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.24'
$ws.Range("E45").Value = '  -4.27%  '
$ws.Range("E46").Value = '  -6.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.54'
$ws.Range("E47").Value = '  -6.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.40'
$ws.Range("E48").Value = '  -10.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.61'
$ws.Range("E49").Value = '  -3.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.856'
$ws.Range("E50").Value = '  -7.32%  '
$ws.Range("D51").Value = '2.194.47'
$ws.Range("E51").Value = '  -8.74%  '
